$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need the
# NumberFormat forced to Text ("@") first, otherwise Excel auto-converts
# the literal price string into a floating point number (losing trailing
# zeros / exact text, e.g. "49.80" -> 49.8).
$textCells = @("D5", "D6", "D10", "D12", "D14", "D18", "D20", "D21", "D23", "D24", "D25", "D27", "D29", "D30", "D33", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D45", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '65.944.19'
$ws.Range("E2").Value = '  +6.48%  '

# Row 3
$ws.Range("D3").Value = '3.010.58'
$ws.Range("E3").Value = '  +3.78%  '

# Row 5
$ws.Range("D5").Value = '583.34'
$ws.Range("E5").Value = '  +2.52%  '

# Row 6
$ws.Range("D6").Value = '162.96'
$ws.Range("E6").Value = '  +13.39%  '

# Row 7
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$ws.Range("E8").Value = '  +3.32%  '

# Row 9
$ws.Range("D9").Value = '3.006.84'
$ws.Range("E9").Value = '  +3.70%  '

# Row 10
$ws.Range("D10").Value = '6.87'
$ws.Range("E10").Value = '  -0.78%  '

# Row 11
$ws.Range("E11").Value = '  +7.65%  '

# Row 12
$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  +6.75%  '

# Row 13
$ws.Range("E13").Value = '  +8.89%  '

# Row 14
$ws.Range("D14").Value = '34.77'
$ws.Range("E14").Value = '  +8.03%  '

# Row 15
$ws.Range("E15").Value = '  -0.47%  '

# Row 16
$ws.Range("D16").Value = '65.921.70'
$ws.Range("E16").Value = '  +6.56%  '

# Row 17
$ws.Range("D17").Value = '3.510.73'
$ws.Range("E17").Value = '  +3.83%  '

# Row 18
$ws.Range("D18").Value = '6.97'
$ws.Range("E18").Value = '  +7.00%  '

# Row 19
$ws.Range("D19").Value = '3.013.93'
$ws.Range("E19").Value = '  +3.96%  '

# Row 20
$ws.Range("D20").Value = '457.77'
$ws.Range("E20").Value = '  +6.36%  '

# Row 21
$ws.Range("D21").Value = '13.96'
$ws.Range("E21").Value = '  +8.12%  '

# Row 22
$ws.Range("E22").Value = '  +5.57%  '

# Row 23
$ws.Range("D23").Value = '7.37'
$ws.Range("E23").Value = '  +7.33%  '

# Row 24
$ws.Range("D24").Value = '82.56'
$ws.Range("E24").Value = '  +4.52%  '

# Row 25
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  +14.06%  '

# Row 26
$ws.Range("E26").Value = '  +3.36%  '

# Row 27
$ws.Range("D27").Value = '10.64'
$ws.Range("E27").Value = '  +5.36%  '

# Row 29
$ws.Range("D29").Value = '8.15'
$ws.Range("E29").Value = '  +16.13%  '

# Row 30
$ws.Range("D30").Value = '2.34'
$ws.Range("E30").Value = '  +15.69%  '

# Row 31
$ws.Range("E31").Value = '  +3.84%  '

# Row 32
$ws.Range("E32").Value = '  -6.74%  '

# Row 33
$ws.Range("D33").Value = '27.06'
$ws.Range("E33").Value = '  +5.57%  '

# Row 34
$ws.Range("E34").Value = '  +3.68%  '

# Row 35
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.05%  '

# Row 36
$ws.Range("D36").Value = '0.994'
$ws.Range("E36").Value = '  +4.05%  '

# Row 37
$ws.Range("D37").Value = '5.81'
$ws.Range("E37").Value = '  +7.65%  '

# Row 38
$ws.Range("E38").Value = '  +11.97%  '

# Row 39
$ws.Range("E39").Value = '  +5.67%  '

# Row 40
$ws.Range("D40").Value = '49.80'
$ws.Range("E40").Value = '  +1.99%  '

# Row 41
$ws.Range("D41").Value = '0.309'
$ws.Range("E41").Value = '  +15.05%  '

# Row 42
$ws.Range("B42").Value = 'Arweave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D42").Value = '43.72'
$ws.Range("E42").Value = '  +7.86%  '

# Row 43
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").Value = '0.121'
$ws.Range("E43").Value = '  +5.67%  '

# Row 44
$ws.Range("E44").Value = '  +4.03%  '

# Row 45
$ws.Range("D45").Value = '388.79'
$ws.Range("E45").Value = '  +12.51%  '

# Row 46
$ws.Range("E46").Value = '  +6.02%  '

# Row 47
$ws.Range("D47").Value = '2.796.55'

# Row 48
$ws.Range("D48").Value = '135.05'
$ws.Range("E48").Value = '  +2.52%  '

# Row 50
$ws.Range("D50").Value = '23.91'
$ws.Range("E50").Value = '  +10.83%  '

# Row 51
$ws.Range("E51").Value = '  +4.03%  '
